$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B231:B325").Value = 1
$ws.Range("B464:B518").Value = 2
